$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new Mac-Address rows (31 and 32), following the same pattern
# as the existing rows: regcntr_id, usr_id, lang_code, is_active, cr_by, cr_dtimes
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 110030
$ws.Range("C31").Value = "eng"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"

$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 110031
$ws.Range("C32").Value = "eng"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"

# Move the active selection to D32, matching where the user ended up
# after entering the new data
$ws.Range("D32").Select()
$excel.ActiveWindow.ScrollRow = 16
